$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 0) Remove the existing hyperlink (old H2 -> tom@gmail.com); we'll re-add
#    hyperlinks for every email cell in the new layout further down.
# ---------------------------------------------------------------------------
foreach ($h in $ws.Hyperlinks) { $h.Delete() }

# ---------------------------------------------------------------------------
# 1) Stash the existing header look (bold font + fill) from A1 onto a scratch
#    cell well outside the sheet's future used range, so we can reuse the
#    exact same style (and not create a brand-new fill/font entry) once the
#    old data has been wiped.
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2) Wipe out all of the old data/format (but keep our Z1 scratch cell).
# ---------------------------------------------------------------------------
$ws.Range("A1:I2").Clear()

# ---------------------------------------------------------------------------
# 3) Re-apply the header style across the new header row A1:K1.
# ---------------------------------------------------------------------------
$ws.Range("Z1").Copy()
$ws.Range("A1:K1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Scratch cell no longer needed.
$ws.Range("Z1").Clear()

# ---------------------------------------------------------------------------
# 4) Header row values.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "firstname"
$ws.Range("B1").Value = "lastname"
$ws.Range("C1").Value = "phone"
$ws.Range("D1").Value = "email"
$ws.Range("E1").Value = "address1"
$ws.Range("F1").Value = "address2"
$ws.Range("G1").Value = "city"
$ws.Range("H1").Value = "state"
$ws.Range("I1").Value = "zipcode"
$ws.Range("J1").Value = "username"
$ws.Range("K1").Value = "password"

# ---------------------------------------------------------------------------
# 5) Data rows.
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "Tom"
$ws.Range("B2").Value = "Rider"
$ws.Range("C2").Value = 9876543210
$ws.Range("D2").Value = "tom@gmail.com"
$ws.Range("E2").Value = "12, high street"
$ws.Range("F2").Value = "downtown"
$ws.Range("G2").Value = "Sunnyvale"
$ws.Range("H2").Value = "California"
$ws.Range("I2").Value = 123456
$ws.Range("J2").Value = "Tom"
$ws.Range("K2").Value = "Rider"

$ws.Range("A3").Value = "Jack"
$ws.Range("B3").Value = "Ram"
$ws.Range("C3").Value = 9876542322
$ws.Range("D3").Value = "jack@gmail.com"
$ws.Range("E3").Value = "12, high street"
$ws.Range("F3").Value = "downtown"
$ws.Range("G3").Value = "Sunnyvale"
$ws.Range("H3").Value = "California"
$ws.Range("I3").Value = 123456
$ws.Range("J3").Value = "Jack"
$ws.Range("K3").Value = "Ram"

$ws.Range("A4").Value = "Jill"
$ws.Range("B4").Value = "jones"
$ws.Range("C4").Value = 9876542223
$ws.Range("D4").Value = "jill@gmail.com"
$ws.Range("E4").Value = "12, high street"
$ws.Range("F4").Value = "downtown"
$ws.Range("G4").Value = "Sunnyvale"
$ws.Range("H4").Value = "California"
$ws.Range("I4").Value = 123456
$ws.Range("J4").Value = "Jill"
$ws.Range("K4").Value = "jones"

$ws.Range("A5").Value = "James"
$ws.Range("B5").Value = "Bond"
$ws.Range("C5").Value = 9876543233
$ws.Range("D5").Value = "james@gmail.com"
$ws.Range("E5").Value = "12, high street"
$ws.Range("F5").Value = "downtown"
$ws.Range("G5").Value = "Sunnyvale"
$ws.Range("H5").Value = "California"
$ws.Range("I5").Value = 123456
$ws.Range("J5").Value = "James"
$ws.Range("K5").Value = "Bond"

# ---------------------------------------------------------------------------
# 6) Hyperlink the email column for every data row.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:tom@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:jack@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:jill@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:james@gmail.com")
$ws.Range("D2:D5").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# 7) Column widths for the re-shuffled columns C, D, E.
# ---------------------------------------------------------------------------
$ws.Columns("C").ColumnWidth = 12.16
$ws.Columns("D").ColumnWidth = 17.29
$ws.Columns("E").ColumnWidth = 15.16

# ---------------------------------------------------------------------------
# 8) Selection ends on K5.
# ---------------------------------------------------------------------------
$ws.Range("K5").Select()
